# Atualizei dados da bibi
# Update retention metrics for a few cohort rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: num_customers 37 -> 38; retention_rate = num_customers / cohort_size
$ws.Range("C27").Value = 38
$ws.Range("E27").Value = 38 / 2252

# Row 36: num_customers 105 -> 106; retention_rate = num_customers / cohort_size
$ws.Range("C36").Value = 106
$ws.Range("E36").Value = 106 / 1930

# Row 37: num_customers 634 -> 648; cohort_size 634 -> 648 (retention_rate stays 1)
$ws.Range("C37").Value = 648
$ws.Range("D37").Value = 648
